$wb = $excel.ActiveWorkbook

# --- Entities sheet: mark A12 (QuestionTimeline) with the "Good" style ---
$wsEntities = $wb.Worksheets.Item("Entities")
$wsEntities.Cells.Item(12, 1).Style = "Good"
[void]$wsEntities.Range("A12").Select()

# --- Methods sheet: mark the finished "question" (and other) method groups as Done ---
$wsMethods = $wb.Worksheets.Item("Methods")
$doneRows = @(6,7,8,9,10,11,21,22,23,24,28,29,30,31,32,33,34,35,36)
foreach ($r in $doneRows) {
    $cell = $wsMethods.Cells.Item($r, 1)
    $cell.Value = "Done"
    $cell.Style = "Good"
}

# Methods becomes the active/selected sheet & tab, with A36 as the active cell
$wsMethods.Activate()
[void]$wsMethods.Range("A36").Select()
